# Generate Report for Handback
#
# The CI handback report is regenerated: the row for
# 10e5e205-878b-46b5-ac70-bf4b698a19c4 moves from "Ready for handoff" to
# "Handed back: in sync with en-US" on every sheet, and the per-language
# sheets (zh-cn, de-de) gain their "Latest Target File" / "Latest Handback
# File" hyperlinks plus an updated "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# Cornflower blue (FF6495ED) underlined font, matching the workbook's
# existing custom "HyperLink" look used by the other linked cells.
$linkColor = 15570276

function Style-AsLink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $linkColor
}

# ---------------------------------------------------------------------
# Overview sheet: the summary status for 10e5e205-...md is now "handed
# back" in both the zh-cn and de-de columns.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusHandedBack
$overview.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusHandedBack

$zhcn.Range("F2").Value = "10e5e205-878b-46b5-ac70-bf4b698a19c4.md"
Style-AsLink $zhcn.Range("F2")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/115ca7257b8688db56f0dda7fc43243b412b4ef6/e2e/10e5e205-878b-46b5-ac70-bf4b698a19c4.md", "", "", "10e5e205-878b-46b5-ac70-bf4b698a19c4.md")

$zhcn.Range("G2").Value = "10e5e205-878b-46b5-ac70-bf4b698a19c4.d7b5908baf43d5545d260898841ea08ef89a6409.zh-cn.xlf"
Style-AsLink $zhcn.Range("G2")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b5e26e155595ce35279441be2f9e4adb61247aa7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/10e5e205-878b-46b5-ac70-bf4b698a19c4.d7b5908baf43d5545d260898841ea08ef89a6409.zh-cn.xlf", "", "", "10e5e205-878b-46b5-ac70-bf4b698a19c4.d7b5908baf43d5545d260898841ea08ef89a6409.zh-cn.xlf")

$zhcn.Range("H2").Value = "2016-03-19 20:29:31"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusHandedBack

$dede.Range("F2").Value = "10e5e205-878b-46b5-ac70-bf4b698a19c4.md"
Style-AsLink $dede.Range("F2")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/115ca7257b8688db56f0dda7fc43243b412b4ef6/e2e/10e5e205-878b-46b5-ac70-bf4b698a19c4.md", "", "", "10e5e205-878b-46b5-ac70-bf4b698a19c4.md")

$dede.Range("G2").Value = "10e5e205-878b-46b5-ac70-bf4b698a19c4.d7b5908baf43d5545d260898841ea08ef89a6409.de-de.xlf"
Style-AsLink $dede.Range("G2")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67d6dc93b07adf47fd05a3b2d4310146f3ce0cb8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/10e5e205-878b-46b5-ac70-bf4b698a19c4.d7b5908baf43d5545d260898841ea08ef89a6409.de-de.xlf", "", "", "10e5e205-878b-46b5-ac70-bf4b698a19c4.d7b5908baf43d5545d260898841ea08ef89a6409.de-de.xlf")

$dede.Range("H2").Value = "2016-03-19 20:29:36"
